# Revert "new changes in ops (ordercreation & orderpage & order form)"
#
# Net effect on the "with_all_correctdata" sheet:
#   - Column "Client" (currently column G) moves to column E, pushing
#     "Typist" (E->F) and "Typist QC" (F->G) one column to the right.
#   - "Product Name" (column J) and "Lob" (column H) swap places.
#   - The selected cell changes from F4 to E12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Use a scratch area far outside the used range (A1:M3) as staging,
#     so Copy/Paste round-trips (which preserve both value AND style)
#     don't clobber each other when we rotate E/F/G. ---
$stageE = $ws.Range("AA1:AA3")
$stageF = $ws.Range("AB1:AB3")
$stageG = $ws.Range("AC1:AC3")
$stageH = $ws.Range("AD1:AD3")

# Stage current E, F, G, H columns (rows 1-3)
$ws.Range("E1:E3").Copy($stageE)
$ws.Range("F1:F3").Copy($stageF)
$ws.Range("G1:G3").Copy($stageG)
$ws.Range("H1:H3").Copy($stageH)

# Rotate: E <- old G, F <- old E, G <- old F  (the "Client" column move)
$stageG.Copy($ws.Range("E1:E3"))
$stageE.Copy($ws.Range("F1:F3"))
$stageF.Copy($ws.Range("G1:G3"))

# Swap H <-> J ("Lob" <-> "Product Name"); stageH still holds the
# original H (pre-rotation H/J are untouched by the E/F/G rotation).
$ws.Range("J1:J3").Copy($ws.Range("H1:H3"))
$stageH.Copy($ws.Range("J1:J3"))

# Clean up the scratch area.
$ws.Range("AA1:AD3").Clear()

# --- Column width: the "Client" column's explicit width (12.6640625)
#     now belongs under column E instead of column G. ---
$ws.Columns(5).ColumnWidth = 11.83

# --- Selection moves from F4 to E12. ---
$ws.Range("E12").Select() | Out-Null
